$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Fail" status values from B2:B53 (header "Status" in B1 stays)
$ws.Range("B2:B53").ClearContents()

# Update the active selection to match the edited workbook state
$ws.Range("J10").Select()
